$wb = $excel.ActiveWorkbook

# The bfbf0c67 file has been handed back (in sync with en-US) for both
# locales. Update status + handback timestamps to reflect the new report.

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the bfbf0c67-*.md file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet: row 3 is the bfbf0c67-*.md file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-03-17 06:32:41"

# --- de-de sheet: row 3 is the bfbf0c67-*.md file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-03-17 06:32:54"
